# Apply the "added selectNext, selectAll, and clearAll" edits to the jobs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("jobs")

# Row 2: Consulting Data Engineer / The Cigna Group
#   end date pushed out a month, and text/RGB colors swapped
$ws.Range("D2").Value = 45139
$ws.Range("G2").Value = "#006688"
$ws.Range("I2").Value = "white"

# Row 4: Consulting Data Engineer / Angel Studios
#   css RGB changed
$ws.Range("G4").Value = "#4400cd"

# Row 7: Senior Data Engineer / SeniorLiink -> b
$ws.Range("B7").Value = "b"

# Re-wrapped row heights (Excel autofit recalculation triggered by the edits)
$ws.Rows.Item(7).RowHeight = 212.25
$ws.Rows.Item(8).RowHeight = 88.5
$ws.Rows.Item(9).RowHeight = 212.25
$ws.Rows.Item(10).RowHeight = 184.5
